$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-26 Thursday", $true, $true, $false, $false, $false, $true, 1, $false, "2024-12-27 Friday", 2) | Out-Null
$d.Content.Find.Execute("66+30=", $true, $true, $false, $false, $false, $true, 1, $false, "67-59=", 2) | Out-Null
$d.Content.Find.Execute("42+51=", $true, $true, $false, $false, $false, $true, 1, $false, "78+21=", 2) | Out-Null
$d.Content.Find.Execute("44+13=", $true, $true, $false, $false, $false, $true, 1, $false, "74+10=", 2) | Out-Null
$d.Content.Find.Execute("27-19=", $true, $true, $false, $false, $false, $true, 1, $false, "35+1=", 2) | Out-Null
$d.Content.Find.Execute("47+34=", $true, $true, $false, $false, $false, $true, 1, $false, "76-51=", 2) | Out-Null
$d.Content.Find.Execute("86-1=", $true, $true, $false, $false, $false, $true, 1, $false, "94-8=", 2) | Out-Null
$d.Content.Find.Execute("70-20=", $true, $true, $false, $false, $false, $true, 1, $false, "42+15=", 2) | Out-Null
$d.Content.Find.Execute("94-83=", $true, $true, $false, $false, $false, $true, 1, $false, "41+27=", 2) | Out-Null
$d.Content.Find.Execute("9+41=", $true, $true, $false, $false, $false, $true, 1, $false, "70-46=", 2) | Out-Null
$d.Content.Find.Execute("17+38=", $true, $true, $false, $false, $false, $true, 1, $false, "44+48=", 2) | Out-Null
$d.Content.Find.Execute("6+41=", $true, $true, $false, $false, $false, $true, 1, $false, "51+27=", 2) | Out-Null
$d.Content.Find.Execute("80-66=", $true, $true, $false, $false, $false, $true, 1, $false, "54+12=", 2) | Out-Null
$d.Content.Find.Execute("76-63=", $true, $true, $false, $false, $false, $true, 1, $false, "87-32=", 2) | Out-Null
$d.Content.Find.Execute("19-9=", $true, $true, $false, $false, $false, $true, 1, $false, "51-34=", 2) | Out-Null
$d.Content.Find.Execute("83-60=", $true, $true, $false, $false, $false, $true, 1, $false, "2+62=", 2) | Out-Null
$d.Content.Find.Execute("88+1=", $true, $true, $false, $false, $false, $true, 1, $false, "19+27=", 2) | Out-Null
$d.Content.Find.Execute("85-71=", $true, $true, $false, $false, $false, $true, 1, $false, "9+10=", 2) | Out-Null
$d.Content.Find.Execute("77+15=", $true, $true, $false, $false, $false, $true, 1, $false, "77-47=", 2) | Out-Null
$d.Content.Find.Execute("66-16=", $true, $true, $false, $false, $false, $true, 1, $false, "55+0=", 2) | Out-Null
$d.Content.Find.Execute("74+3=", $true, $true, $false, $false, $false, $true, 1, $false, "72-54=", 2) | Out-Null
$d.Content.Find.Execute("91-88=", $true, $true, $false, $false, $false, $true, 1, $false, "4+65=", 2) | Out-Null
$d.Content.Find.Execute("38+5=", $true, $true, $false, $false, $false, $true, 1, $false, "97-31=", 2) | Out-Null
$d.Content.Find.Execute("16+81=", $true, $true, $false, $false, $false, $true, 1, $false, "41+36=", 2) | Out-Null
$d.Content.Find.Execute("72-4=", $true, $true, $false, $false, $false, $true, 1, $false, "77+18=", 2) | Out-Null
$d.Content.Find.Execute("60-18=", $true, $true, $false, $false, $false, $true, 1, $false, "30-20=", 2) | Out-Null
$d.Content.Find.Execute("31+7=", $true, $true, $false, $false, $false, $true, 1, $false, "78+9=", 2) | Out-Null
$d.Content.Find.Execute("31+23=", $true, $true, $false, $false, $false, $true, 1, $false, "18-16=", 2) | Out-Null
$d.Content.Find.Execute("47+36=", $true, $true, $false, $false, $false, $true, 1, $false, "69+26=", 2) | Out-Null
$d.Content.Find.Execute("9+27=", $true, $true, $false, $false, $false, $true, 1, $false, "77-53=", 2) | Out-Null
$d.Content.Find.Execute("12+33=", $true, $true, $false, $false, $false, $true, 1, $false, "82+12=", 2) | Out-Null
$d.Content.Find.Execute("60-28=", $true, $true, $false, $false, $false, $true, 1, $false, "95-34=", 2) | Out-Null
$d.Content.Find.Execute("5+39=", $true, $true, $false, $false, $false, $true, 1, $false, "25+74=", 2) | Out-Null
$d.Content.Find.Execute("87-9=", $true, $true, $false, $false, $false, $true, 1, $false, "58-50=", 2) | Out-Null
$d.Content.Find.Execute("59-5=", $true, $true, $false, $false, $false, $true, 1, $false, "93-22=", 2) | Out-Null
$d.Content.Find.Execute("70-62=", $true, $true, $false, $false, $false, $true, 1, $false, "11+0=", 2) | Out-Null
$d.Content.Find.Execute("74+14=", $true, $true, $false, $false, $false, $true, 1, $false, "17+61=", 2) | Out-Null
$d.Content.Find.Execute("92-57=", $true, $true, $false, $false, $false, $true, 1, $false, "84-19=", 2) | Out-Null
$d.Content.Find.Execute("55+5=", $true, $true, $false, $false, $false, $true, 1, $false, "98-30=", 2) | Out-Null
$d.Content.Find.Execute("26+59=", $true, $true, $false, $false, $false, $true, 1, $false, "71+7=", 2) | Out-Null
$d.Content.Find.Execute("40-0=", $true, $true, $false, $false, $false, $true, 1, $false, "8+65=", 2) | Out-Null
$d.Content.Find.Execute("34+31=", $true, $true, $false, $false, $false, $true, 1, $false, "82-13=", 2) | Out-Null
$d.Content.Find.Execute("89+4=", $true, $true, $false, $false, $false, $true, 1, $false, "71-48=", 2) | Out-Null
$d.Content.Find.Execute("33-19=", $true, $true, $false, $false, $false, $true, 1, $false, "28+60=", 2) | Out-Null
$d.Content.Find.Execute("0+65=", $true, $true, $false, $false, $false, $true, 1, $false, "65-52=", 2) | Out-Null
$d.Content.Find.Execute("78-22=", $true, $true, $false, $false, $false, $true, 1, $false, "70+6=", 2) | Out-Null
$d.Content.Find.Execute("6+80=", $true, $true, $false, $false, $false, $true, 1, $false, "65+18=", 2) | Out-Null
$d.Content.Find.Execute("72+2=", $true, $true, $false, $false, $false, $true, 1, $false, "27+11=", 2) | Out-Null
$d.Content.Find.Execute("81-23=", $true, $true, $false, $false, $false, $true, 1, $false, "44-0=", 2) | Out-Null
$d.Content.Find.Execute("55-45=", $true, $true, $false, $false, $false, $true, 1, $false, "11-3=", 2) | Out-Null
$d.Content.Find.Execute("69+2=", $true, $true, $false, $false, $false, $true, 1, $false, "90-85=", 2) | Out-Null
$d.Content.Find.Execute("83-33=", $true, $true, $false, $false, $false, $true, 1, $false, "4+50=", 2) | Out-Null
$d.Content.Find.Execute("54+1=", $true, $true, $false, $false, $false, $true, 1, $false, "34-13=", 2) | Out-Null
$d.Content.Find.Execute("4+79=", $true, $true, $false, $false, $false, $true, 1, $false, "93-43=", 2) | Out-Null
$d.Content.Find.Execute("56-8=", $true, $true, $false, $false, $false, $true, 1, $false, "96-71=", 2) | Out-Null
$d.Content.Find.Execute("77-11=", $true, $true, $false, $false, $false, $true, 1, $false, "5+9=", 2) | Out-Null
$d.Content.Find.Execute("21-10=", $true, $true, $false, $false, $false, $true, 1, $false, "83-11=", 2) | Out-Null
$d.Content.Find.Execute("88-54=", $true, $true, $false, $false, $false, $true, 1, $false, "2+1=", 2) | Out-Null
$d.Content.Find.Execute("47+47=", $true, $true, $false, $false, $false, $true, 1, $false, "65+28=", 2) | Out-Null
$d.Content.Find.Execute("12+48=", $true, $true, $false, $false, $false, $true, 1, $false, "75-53=", 2) | Out-Null
$d.Content.Find.Execute("54+15=", $true, $true, $false, $false, $false, $true, 1, $false, "32+65=", 2) | Out-Null
$d.Content.Find.Execute("40+42=", $true, $true, $false, $false, $false, $true, 1, $false, "2+33=", 2) | Out-Null
$d.Content.Find.Execute("84-27=", $true, $true, $false, $false, $false, $true, 1, $false, "10+83=", 2) | Out-Null
$d.Content.Find.Execute("56+22=", $true, $true, $false, $false, $false, $true, 1, $false, "92-18=", 2) | Out-Null
$d.Content.Find.Execute("51+30=", $true, $true, $false, $false, $false, $true, 1, $false, "29+13=", 2) | Out-Null
$d.Content.Find.Execute("19+66=", $true, $true, $false, $false, $false, $true, 1, $false, "60+14=", 2) | Out-Null
$d.Content.Find.Execute("27-11=", $true, $true, $false, $false, $false, $true, 1, $false, "8+56=", 2) | Out-Null
$d.Content.Find.Execute("19+73=", $true, $true, $false, $false, $false, $true, 1, $false, "96-83=", 2) | Out-Null
$d.Content.Find.Execute("14+81=", $true, $true, $false, $false, $false, $true, 1, $false, "47+42=", 2) | Out-Null
$d.Content.Find.Execute("4+17=", $true, $true, $false, $false, $false, $true, 1, $false, "10+1=", 2) | Out-Null
$d.Content.Find.Execute("87+6=", $true, $true, $false, $false, $false, $true, 1, $false, "74-1=", 2) | Out-Null
$d.Content.Find.Execute("90-3=", $true, $true, $false, $false, $false, $true, 1, $false, "96-54=", 2) | Out-Null
$d.Content.Find.Execute("29+11=", $true, $true, $false, $false, $false, $true, 1, $false, "8+15=", 2) | Out-Null
$d.Content.Find.Execute("40+40=", $true, $true, $false, $false, $false, $true, 1, $false, "97-27=", 2) | Out-Null
$d.Content.Find.Execute("37+40=", $true, $true, $false, $false, $false, $true, 1, $false, "39+59=", 2) | Out-Null
$d.Content.Find.Execute("95-1=", $true, $true, $false, $false, $false, $true, 1, $false, "77-38=", 2) | Out-Null
$d.Content.Find.Execute("42-32=", $true, $true, $false, $false, $false, $true, 1, $false, "84-78=", 2) | Out-Null
$d.Content.Find.Execute("84+7=", $true, $true, $false, $false, $false, $true, 1, $false, "38+34=", 2) | Out-Null
$d.Content.Find.Execute("99-41=", $true, $true, $false, $false, $false, $true, 1, $false, "7+34=", 2) | Out-Null
$d.Content.Find.Execute("28-16=", $true, $true, $false, $false, $false, $true, 1, $false, "32+26=", 2) | Out-Null
$d.Content.Find.Execute("29-27=", $true, $true, $false, $false, $false, $true, 1, $false, "35+53=", 2) | Out-Null
$d.Content.Find.Execute("87-84=", $true, $true, $false, $false, $false, $true, 1, $false, "34+24=", 2) | Out-Null
$d.Content.Find.Execute("75-37=", $true, $true, $false, $false, $false, $true, 1, $false, "53-21=", 2) | Out-Null
$d.Content.Find.Execute("46+41=", $true, $true, $false, $false, $false, $true, 1, $false, "49+28=", 2) | Out-Null
$d.Content.Find.Execute("6+79=", $true, $true, $false, $false, $false, $true, 1, $false, "96-11=", 2) | Out-Null
$d.Content.Find.Execute("45-27=", $true, $true, $false, $false, $false, $true, 1, $false, "45+35=", 2) | Out-Null
$d.Content.Find.Execute("39+39=", $true, $true, $false, $false, $false, $true, 1, $false, "43+28=", 2) | Out-Null
$d.Content.Find.Execute("77+3=", $true, $true, $false, $false, $false, $true, 1, $false, "6-4=", 2) | Out-Null
$d.Content.Find.Execute("0+41=", $true, $true, $false, $false, $false, $true, 1, $false, "48-5=", 2) | Out-Null
$d.Content.Find.Execute("39-5=", $true, $true, $false, $false, $false, $true, 1, $false, "36+8=", 2) | Out-Null
$d.Content.Find.Execute("37+52=", $true, $true, $false, $false, $false, $true, 1, $false, "1+33=", 2) | Out-Null
$d.Content.Find.Execute("23+59=", $true, $true, $false, $false, $false, $true, 1, $false, "63-52=", 2) | Out-Null
$d.Content.Find.Execute("55+33=", $true, $true, $false, $false, $false, $true, 1, $false, "68+8=", 2) | Out-Null
$d.Content.Find.Execute("37+1=", $true, $true, $false, $false, $false, $true, 1, $false, "81-59=", 2) | Out-Null
$d.Content.Find.Execute("76-14=", $true, $true, $false, $false, $false, $true, 1, $false, "7+27=", 2) | Out-Null
$d.Content.Find.Execute("25+37=", $true, $true, $false, $false, $false, $true, 1, $false, "64+14=", 2) | Out-Null
$d.Content.Find.Execute("86-31=", $true, $true, $false, $false, $false, $true, 1, $false, "78-74=", 2) | Out-Null
$d.Content.Find.Execute("50-26=", $true, $true, $false, $false, $false, $true, 1, $false, "83+2=", 2) | Out-Null
$d.Content.Find.Execute("16+72=", $true, $true, $false, $false, $false, $true, 1, $false, "53-36=", 2) | Out-Null
$d.Content.Find.Execute("5+0=", $true, $true, $false, $false, $false, $true, 1, $false, "19+60=", 2) | Out-Null
$d.Content.Find.Execute("9+74=", $true, $true, $false, $false, $false, $true, 1, $false, "32+2=", 2) | Out-Null
